$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.751.40"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "1.745.46"
$ws.Range("E3").Value = "  -4.44%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "320.15"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "0.4189"
$ws.Range("E7").Value = "  -5.86%  "
$ws.Range("D8").Value = "0.3625"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").Value = "42.72"
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").Value = "0.07409"
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "1.085"
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "20.62"
$ws.Range("E13").Value = "  -7.10%  "
$ws.Range("D14").Value = "6.026"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").Value = "7.257"
$ws.Range("E15").Value = "  -3.95%  "
$ws.Range("D16").Value = "1.768.61"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "90.53"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").Value = "0.00001047"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").Value = "0.06326"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "16.96"
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("D22").Value = "5.913"
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("D23").Value = "27.765.35"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  -5.07%  "
$ws.Range("D25").Value = "2.089"
$ws.Range("E25").Value = "  -3.62%  "
$ws.Range("D26").Value = "156.87"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").Value = "1.965.44"
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("D29").Value = "2.128"
$ws.Range("E29").Value = "  -9.72%  "
$ws.Range("D30").Value = "123.30"
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("D31").Value = "1.123"
$ws.Range("E31").Value = "  -6.08%  "
$ws.Range("D32").Value = "3.638"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "5.535"
$ws.Range("E33").Value = "  -6.15%  "
$ws.Range("D34").Value = "0.08821"
$ws.Range("E34").Value = "  -4.98%  "
$ws.Range("D35").Value = "12.26"
$ws.Range("E35").Value = "  -6.72%  "
$ws.Range("D36").Value = "0.02266"
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("D37").Value = "0.2088"
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").Value = "0.05981"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").Value = "4.938"
$ws.Range("E39").Value = "  -4.76%  "
$ws.Range("D40").Value = "0.6273"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").Value = "1.169"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "1.393"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "7.783"
$ws.Range("E44").Value = "  -4.25%  "
$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  -5.47%  "
$ws.Range("D46").Value = "0.5845"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").Value = "3.666"
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "121.90"
$ws.Range("E48").Value = "  -4.18%  "
$ws.Range("D49").Value = "1.960"
$ws.Range("E49").Value = "  -4.45%  "
$ws.Range("D50").Value = "1.172"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "0.06790"
$ws.Range("E51").Value = "  -3.01%  "
